$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column (C) for rows 2-12
# from 45233 (2023-11-03) to 45243 (2023-11-13)
$ws.Range("C2:C12").Value = 45243
